$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-12 (values change, country names stay the same)
$ws.Range("B2").Value = 68
$ws.Range("C2").Value = 23.53

$ws.Range("B3").Value = 53
$ws.Range("C3").Value = 18.34

$ws.Range("A4").Value = "France"
$ws.Range("B4").Value = 29
$ws.Range("C4").Value = 10.03

$ws.Range("A5").Value = "Sweden"
$ws.Range("B5").Value = 29
$ws.Range("C5").Value = 10.03

$ws.Range("B6").Value = 24
$ws.Range("C6").Value = 8.300000000000001

$ws.Range("B7").Value = 18
$ws.Range("C7").Value = 6.23

$ws.Range("B8").Value = 14
$ws.Range("C8").Value = 4.84

$ws.Range("B9").Value = 13
$ws.Range("C9").Value = 4.5

$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 3.46

$ws.Range("B11").Value = 7
$ws.Range("C11").Value = 2.42

$ws.Range("B12").Value = 6
$ws.Range("C12").Value = 2.08

$ws.Range("A13").Value = "Denmark"
$ws.Range("B13").Value = 6
$ws.Range("C13").Value = 2.08

$ws.Range("A14").Value = "Norway"
$ws.Range("B14").Value = 4
$ws.Range("C14").Value = 1.38

$ws.Range("A15").Value = "Austria"
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 1.38

$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 0.6899999999999999

# New row 17
$ws.Range("A17").Value = "Greece"
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = 0.6899999999999999

# Copy style from A16 to A17 to match formatting (bold, border, centered)
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)  # xlPasteFormats
